# "Activated PageProcessor. Added SCRIPT capabilities (column)."
#
# Target sheet is Sheet3 ("PageProcessor") - make/keep it the active sheet,
# add a SCRIPT column (inserted before the existing BACK_URL column), fix a
# couple of existing values, append a new row-4 (search-results) scenario,
# and turn the rich-text e-mail hyperlink cell into a plain text + hyperlink
# cell with an updated display string.

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()

# --- Insert the new "SCRIPT" column before the old BACK_URL column ---------
# (BACK_URL shifts from F to G)
$ws3.Columns.Item(6).Insert()
$ws3.Range("F1").Value = "SCRIPT"
$ws3.Range("F2").Value = "try { logs.info(" + [char]0x201C + "\n\n HI \n\n" + [char]0x201D + ");} "

# --- Row 2: single "!" instead of double "!!" in the selector --------------
$ws3.Range("C2").Value = "<@!{xpath=//*[contains(text(),'Hello. Sign in')]}>"

# --- Row 3: RUN flag changed from y -> n ------------------------------------
$ws3.Range("A3").Value = "n"

# --- New row 4: additional "search results" scenario -----------------------
$ws3.Range("A4").Value = "n"
$ws3.Range("C4").Value = "<@!{xpath=//*[@id='twotabsearchtextbox']}>|hair {SPACE} brash"
$ws3.Range("D4").Value = "xpath=//*[@id=" + [char]0x2019 + "didYouMean" + [char]0x2019 + "]"
$ws3.Range("E4").Value = "Showing results for"

# --- C3: replace the two-tone rich text with plain text, refresh the -------
# hyperlink so its display string matches the new value, keep the original
# (non-underlined, blue) font look and the text number format.
$ws3.Hyperlinks.Delete()
$ws3.Range("C3").Value = "gpawel17@mail.com|1qazxsw2!"
$ws3.Hyperlinks.Add($ws3.Range("C3"), "mailto:gpawel17@mail.com", "", "", "gpawel17@mail.com|1qazxsw2!")
$ws3.Range("C3").Font.Underline = $false
$ws3.Range("C3").Font.Name = "Arial"
$ws3.Range("C3").Font.Size = 10
$ws3.Range("C3").NumberFormat = "@"
$ws3.Range("C3").Font.Color = 16711680

# --- Column C is a text column everywhere else, too -------------------------
$ws3.Range("C1").NumberFormat = "@"
$ws3.Range("C2").NumberFormat = "@"
$ws3.Range("C4").NumberFormat = "@"

# --- Restore the selection shown in the saved file --------------------------
$ws3.Range("F7").Select() | Out-Null
